# edit.ps1 - applies the "Ve(This is a change) / Version for main branch" edits
# described by the commit "changes made in main file":
#   1. Append "  (This is a change – Version for main branch)" (the parenthetical
#      in red) to the end of the first paragraph's text.
#   2. Delete the trailing "ank God almighty, we are free at last." paragraph
#      (the one following "Shall be lifted-nevermore!").
#   3. Drop the style definitions that become unused once that paragraph (and
#      its NormalWeb-styled content) disappears.

$d = $word.ActiveDocument

# --- 1. Extend the first paragraph -----------------------------------------
$p1 = $d.Paragraphs.Item(1)
$r = $p1.Range
# Move to just before the paragraph mark so we append inside paragraph 1.
$r.End = $r.End - 1
$r.Collapse(0)

# Two trailing spaces after the existing sentence (plain formatting).
$r.InsertAfter("  ")
$r.Collapse(0)

# The red parenthetical, inserted (and colored) in three separate chunks so
# it ends up as three runs, matching how it was originally typed/edited.
$part1 = "(This is a change " + [char]0x2013 + " Ve"
$part2 = "rsion for main branch"
$part3 = ")"

$start1 = $r.Start
$r.InsertAfter($part1)
$end1 = $start1 + $part1.Length
$d.Range($start1, $end1).Font.Color = 255   # wdColorRed / FF0000

$r.Collapse(0)
$start2 = $r.Start
$r.InsertAfter($part2)
$end2 = $start2 + $part2.Length
$d.Range($start2, $end2).Font.Color = 255

$r.Collapse(0)
$start3 = $r.Start
$r.InsertAfter($part3)
$end3 = $start3 + $part3.Length
$d.Range($start3, $end3).Font.Color = 255

# --- 2. Remove the trailing "...God almighty, we are free at last." paragraph
$last = $d.Paragraphs.Last
$last.Range.Delete()

# --- 3. Remove now-unused style definitions ---------------------------------
# (deleted from the end of the styles collection backwards - deleting forward
# through linked styles, e.g. Heading2/Heading2Char, trips an indexing bug)
$stylesToDelete = @(
  "podcast-toolssubscribe-links",
  "generic-title",
  "subscribe-more-info",
  "subscribe",
  "audio-tool",
  "Heading4Char",
  "Heading2Char",
  "Hyperlink",
  "apple-converted-space",
  "Heading4",
  "Heading2"
)
foreach ($name in $stylesToDelete) {
  $d.Styles.Item($name).Delete()
}
